# Rename the "Growth" worksheet to "Growth_Cold" (ARCtrl reference update),
# then make it the active sheet with the selection parked on G43:H43 - matching
# the state the workbook was left in when it was last saved from Excel.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Growth")
$ws.Name = "Growth_Cold"

$ws.Activate()
$ws.Range("G43:H43").Select()
